# ---------------------------------------------------------------------------
# "feat: add 2022-Q3 data"
#   1. Insert a new "2022-Q3" row into the "总计" (summary) sheet.
#   2. Insert a brand-new worksheet named "2022-Q3" (right after "总计",
#      before "2022-Q2") holding the per-fund holdings detail, matching the
#      layout used by the other quarterly sheets.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ----- locate the sheets we need by name (robust to pre-existing order) ----
$summary = $wb.Worksheets.Item("总计")
$q2      = $wb.Worksheets.Item("2022-Q2")
$lastSheetName = $wb.Worksheets.Item($wb.Worksheets.Count).Name

# =============================================================================
# 1. "总计" sheet: insert the new 2022-Q3 row just under the header row, and
#    push everything else down by one row.
# =============================================================================
$summary.Rows.Item(2).Insert()

# Excel's row-insert copies formatting down from the row above (the header),
# which would leave stray bold/border styling on the new data cells - strip
# it back to the default "no style" used by every other data row.
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 11
$summary.Range("D2").Value = 0.47

# copy the style of the (still intact) index column / row so the new row
# matches the rest of the table
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# renumber the index column (A) for the rows that shifted down, so it stays
# a plain 0..4 sequence
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# =============================================================================
# 2. Brand-new "2022-Q3" worksheet, inserted right before "2022-Q2"
# =============================================================================
$newSheet = $wb.Worksheets.Add($q2)
$newSheet.Name = "2022-Q3"

# ----- header row (bold, boxed, centered - same look as the other sheets) --
$hdr = $newSheet.Range("B1:H1")
$hdr.Value = "基金代码"
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# ----- fund rows -------------------------------------------------------------
$rows = @(
    @("014016", "中信建投品质优选一年持有期混合A", "4.28", "85.84", "2.68", "0.1147", 5),
    @("008347", "中信建投价值甄选混合A",           "3.87", "76.99", "2.44", "0.0944", 6),
    @("003822", "中信建投行业轮换混合A",           "3.10", "77.97", "2.58", "0.0800", 5),
    @("008348", "中信建投价值甄选混合C",           "1.59", "76.99", "2.44", "0.0388", 6),
    @("001097", "华泰柏瑞积极优选股票A",           "1.20", "83.01", "2.73", "0.0328", 10),
    @("014017", "中信建投品质优选一年持有期混合C", "1.20", "85.84", "2.68", "0.0322", 5),
    @("690003", "民生加银精选混合",                 "0.51", "89.96", "5.38", "0.0274", 8),
    @("003823", "中信建投行业轮换混合C",           "0.85", "77.97", "2.58", "0.0219", 5),
    @("007468", "中信建投策略精选混合A",           "0.55", "81.86", "2.59", "0.0142", 4),
    @("007469", "中信建投策略精选混合C",           "0.38", "81.86", "2.59", "0.0098", 4),
    @("016283", "华泰柏瑞积极优选股票C",           "0.16", "83.01", "2.73", "0.0044", 10)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = ($r - 2)
    # columns B-G are stored as text (not numbers) in every other quarter
    # sheet, including fund codes with significant leading zeros - use the
    # leading-apostrophe trick to force text for them.
    $newSheet.Range("B$r").Value = "'" + $row[0]
    $newSheet.Range("C$r").Value = $row[1]
    $newSheet.Range("D$r").Value = "'" + $row[2]
    $newSheet.Range("E$r").Value = "'" + $row[3]
    $newSheet.Range("F$r").Value = "'" + $row[4]
    $newSheet.Range("G$r").Value = "'" + $row[5]
    $newSheet.Range("H$r").Value = $row[6]
    $r = $r + 1
}
$lastRow = $r - 1

# the apostrophe-prefix trick above marks the cells with a stray
# "quote prefix" style - clear it so the data cells end up with the same
# "no explicit style" look used throughout the rest of the workbook.
$newSheet.Range("B2:G$lastRow").ClearFormats()

# index column (A): bold, boxed, centered - same look as the other sheets
$idx = $newSheet.Range("A2:A$lastRow")
$idx.Font.Bold = $true
$idx.Borders.LineStyle = 1
$idx.HorizontalAlignment = -4108
$idx.VerticalAlignment = -4160

# ----- keep the workbook's originally-selected tab selected -----------------
$wb.Worksheets.Item($lastSheetName).Activate()
